$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain stored as text
# (matches the source inlineStr cells in the workbook, e.g. "584.83").
# Setting the NumberFormat to Text ("@") before assigning the value
# prevents Excel from auto-converting the literal into a number.
$textCells = @("D5", "D6", "D8", "D9", "D13", "D18", "D19", "D20", "D21", "D23", "D28", "D29", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.366.18'
$ws.Range("E2").Value = '  +4.15%  '
$ws.Range("D3").Value = '3.502.33'
$ws.Range("E3").Value = '  +4.09%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '584.83'
$ws.Range("D6").Value = '148.03'
$ws.Range("E6").Value = '  +6.96%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.478'
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("D9").Value = '7.68'
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("E10").Value = '  +5.04%  '
$ws.Range("E11").Value = '  +5.29%  '
$ws.Range("D12").Value = '4.098.28'
$ws.Range("E12").Value = '  +3.97%  '
$ws.Range("D13").Value = '29.64'
$ws.Range("E13").Value = '  +7.75%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '3.495.56'
$ws.Range("E15").Value = '  +3.79%  '
$ws.Range("E16").Value = '  +5.09%  '
$ws.Range("D17").Value = '63.454.73'
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").Value = '6.29'
$ws.Range("E18").Value = '  +4.24%  '
$ws.Range("D19").Value = '14.34'
$ws.Range("E19").Value = '  +6.17%  '
$ws.Range("D20").Value = '9.46'
$ws.Range("E20").Value = '  +7.45%  '
$ws.Range("D21").Value = '395.71'
$ws.Range("E21").Value = '  +4.19%  '
$ws.Range("E22").Value = '  +3.70%  '
$ws.Range("D23").Value = '75.42'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +9.38%  '
$ws.Range("D26").Value = '3.643.54'
$ws.Range("E26").Value = '  +3.86%  '
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  +10.10%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '8.29'
$ws.Range("E30").Value = '  +6.10%  '
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("D32").Value = '1.43'
$ws.Range("E32").Value = '  +7.47%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '23.91'
$ws.Range("E34").Value = '  +5.13%  '
$ws.Range("B35").Value = 'EnergySwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D35").Value = '32.69'
$ws.Range("E35").Value = '  +30.67%  '
$ws.Range("D36").Value = '7.21'
$ws.Range("E36").Value = '  +5.63%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '5.37'
$ws.Range("E37").Value = '  +9.87%  '
$ws.Range("D38").Value = '173.15'
$ws.Range("E38").Value = '  +4.37%  '
$ws.Range("E39").Value = '  +10.47%  '
$ws.Range("D40").Value = '3.534.23'
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("D41").Value = '0.0775'
$ws.Range("E41").Value = '  +2.59%  '
$ws.Range("D42").Value = '0.805'
$ws.Range("E42").Value = '  +4.62%  '
$ws.Range("D43").Value = '1.75'
$ws.Range("E43").Value = '  +8.52%  '
$ws.Range("E44").Value = '  +5.25%  '
$ws.Range("D45").Value = '42.29'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  +10.73%  '
$ws.Range("D47").Value = '2.581.69'
$ws.Range("E47").Value = '  +5.70%  '
$ws.Range("D48").Value = '24.06'
$ws.Range("E48").Value = '  +8.94%  '
$ws.Range("D49").Value = '2.27'
$ws.Range("E49").Value = '  +12.80%  '
$ws.Range("D50").Value = '6.78'
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("E51").Value = '  +5.83%  '
